$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the CPU sample log (rows 2-34) with the latest capture pass timestamps/values.
# This adds a flag/window so the next analysis waits for fresh data instead of reusing old rows.

$ws.Range("A2").Value = "05/05/2021 01:39:46"
$ws.Range("B2").Value = 0.5
$ws.Range("A3").Value = "05/05/2021 01:40:48"
$ws.Range("B3").Value = 0.3333
$ws.Range("A4").Value = "05/05/2021 01:41:50"
$ws.Range("B4").Value = 0.3333
$ws.Range("A5").Value = "05/05/2021 01:42:55"
$ws.Range("A6").Value = "05/05/2021 01:44:00"
$ws.Range("B6").Value = 0.3279
$ws.Range("A7").Value = "05/05/2021 01:45:04"
$ws.Range("B7").Value = 43.5593
$ws.Range("A8").Value = "05/05/2021 01:46:09"
$ws.Range("B8").Value = 72.66670000000001
$ws.Range("A9").Value = "05/05/2021 01:51:31"
$ws.Range("B9").Value = 51.01690000000001
$ws.Range("A10").Value = "05/05/2021 01:52:37"
$ws.Range("B10").Value = 30.83335
$ws.Range("A11").Value = "05/05/2021 01:53:43"
$ws.Range("B11").Value = 54.9863
$ws.Range("A12").Value = "05/05/2021 01:54:50"
$ws.Range("B12").Value = 22.13115
$ws.Range("A13").Value = "05/05/2021 02:00:08"
$ws.Range("B13").Value = 68.66670000000001
$ws.Range("A14").Value = "05/05/2021 02:01:17"
$ws.Range("B14").Value = 89.16670000000001
$ws.Range("A15").Value = "05/05/2021 02:06:33"
$ws.Range("B15").Value = 30.25
$ws.Range("A16").Value = "05/05/2021 02:07:38"
$ws.Range("B16").Value = 0.3333
$ws.Range("A17").Value = "05/05/2021 02:12:49"
$ws.Range("B17").Value = 0.5
$ws.Range("A18").Value = "05/05/2021 02:13:53"
$ws.Range("B18").Value = 0.3390000000000001
$ws.Range("A19").Value = "05/05/2021 02:14:57"
$ws.Range("B19").Value = 0.3390000000000001
$ws.Range("A20").Value = "05/05/2021 02:16:00"
$ws.Range("A21").Value = "05/05/2021 02:17:04"
$ws.Range("B21").Value = 0.3333
$ws.Range("A22").Value = "05/05/2021 02:18:07"
$ws.Range("B22").Value = 0.4918
$ws.Range("A23").Value = "05/05/2021 02:19:12"
$ws.Range("B23").Value = 0.3390000000000001
$ws.Range("A24").Value = "05/05/2021 02:20:16"
$ws.Range("B24").Value = 30.50850000000001
$ws.Range("A25").Value = "05/05/2021 02:21:21"
$ws.Range("B25").Value = 38.5
$ws.Range("A26").Value = "05/05/2021 02:21:31"
$ws.Range("B26").Value = 38.5
$ws.Range("A27").Value = "05/05/2021 02:21:40"
$ws.Range("B27").Value = 38.5
$ws.Range("A28").Value = "05/05/2021 02:21:48"
$ws.Range("B28").Value = 38.5
$ws.Range("A29").Value = "05/05/2021 02:21:58"
$ws.Range("B29").Value = 59.3443
$ws.Range("A30").Value = "05/05/2021 02:22:08"
$ws.Range("B30").Value = 59.3443
$ws.Range("A31").Value = "05/05/2021 02:22:18"
$ws.Range("B31").Value = 59.3443
$ws.Range("A32").Value = "05/05/2021 02:22:30"
$ws.Range("B32").Value = 59.3443
$ws.Range("A33").Value = "05/05/2021 02:22:40"
$ws.Range("B33").Value = 59.3443
$ws.Range("A34").Value = "05/05/2021 02:22:49"
$ws.Range("B34").Value = 83.5

# Drop the stale trailing rows from the previous run (old rows 35-53) so data ends at row 34.
$ws.Range("A35:B53").EntireRow.Delete()
